$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 57579.82
$ws.Range("J17").Value = 59081.92
$ws.Range("L17").Value = 177245.76
$ws.Range("N17").Value = -177581.76
$ws.Range("H32").Value = 579.7143
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 509.66666
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 509.66666
$ws.Range("M32").Value = -674
$ws.Range("N32").Value = -1161.66666
$ws.Range("H62").Value = 2695.2307
$ws.Range("I62").Value = 2102.2
$ws.Range("J62").Value = 3065.875
$ws.Range("K62").Value = 2102.2
$ws.Range("L62").Value = 3065.875
$ws.Range("M62").Value = -1478.2
$ws.Range("N62").Value = -4313.875
$ws.Range("H65").Value = 2695.2307
$ws.Range("I65").Value = 2102.2
$ws.Range("J65").Value = 3065.875
$ws.Range("K65").Value = 10511
$ws.Range("L65").Value = 15329.375
$ws.Range("M65").Value = -7391
$ws.Range("N65").Value = -21569.375
$ws.Range("H81").Value = 29206.25
$ws.Range("J81").Value = 29206.25
$ws.Range("L81").Value = 29206.25
$ws.Range("N81").Value = -31202.25
$ws.Range("H84").Value = 29206.25
$ws.Range("J84").Value = 29206.25
$ws.Range("L84").Value = 87618.75
$ws.Range("N84").Value = -97602.75
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H116").Value = 3938.96
$ws.Range("I116").Value = 3298
$ws.Range("J116").Value = 4633.3335
$ws.Range("K116").Value = 3298
$ws.Range("L116").Value = 4633.3335
$ws.Range("M116").Value = 144
$ws.Range("N116").Value = -11517.3335
$ws.Range("H134").Value = 25878.334
$ws.Range("J134").Value = 25878.334
$ws.Range("L134").Value = 25878.334
$ws.Range("N134").Value = -36018.334
$ws.Range("H137").Value = 3765.925
$ws.Range("I137").Value = 4061.926
$ws.Range("J137").Value = 3151.1538
$ws.Range("K137").Value = 12185.778
$ws.Range("L137").Value = 9453.4614
$ws.Range("M137").Value = -9635.778
$ws.Range("N137").Value = -14553.4614
$ws.Range("H141").Value = 416034.03
$ws.Range("I141").Value = 1112.3572
$ws.Range("J141").Value = 692648.5
$ws.Range("K141").Value = 3337.0716
$ws.Range("L141").Value = 2077945.5
$ws.Range("M141").Value = 1842.9284
$ws.Range("N141").Value = -2088305.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2868.7856
$ws.Range("I74").Value = 2021.4445
$ws.Range("K74").Value = 2021.4445
$ws.Range("M74").Value = -1147.4445
$ws.Range("H77").Value = 2868.7856
$ws.Range("I77").Value = 2021.4445
$ws.Range("K77").Value = 10107.2225
$ws.Range("M77").Value = -5739.2225
$ws.Range("H123").Value = 28668.777
$ws.Range("J123").Value = 28668.777
$ws.Range("L123").Value = 28668.777
$ws.Range("N123").Value = -38468.777
$ws.Range("H132").Value = 2836.889
$ws.Range("I132").Value = 2373.2
$ws.Range("J132").Value = 3890.7273
$ws.Range("K132").Value = 7119.599999999999
$ws.Range("L132").Value = 11672.1819
$ws.Range("M132").Value = -4589.599999999999
$ws.Range("N132").Value = -16732.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 1814.3334
$ws.Range("I17").Value = 850
$ws.Range("K17").Value = 850
$ws.Range("M17").Value = -678
$ws.Range("H69").Value = 29431.666
$ws.Range("J69").Value = 29431.666
$ws.Range("L69").Value = 29431.666
$ws.Range("N69").Value = -31053.666
$ws.Range("H72").Value = 29431.666
$ws.Range("J72").Value = 29431.666
$ws.Range("L72").Value = 88294.99800000001
$ws.Range("N72").Value = -96406.99800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3392.4404
$ws.Range("J31").Value = 3586.049
$ws.Range("L31").Value = 3586.049
$ws.Range("N31").Value = -4176.049
$ws.Range("H34").Value = 3392.4404
$ws.Range("J34").Value = 3586.049
$ws.Range("L34").Value = 3586.049
$ws.Range("N34").Value = -3990.049
$ws.Range("H68").Value = 29929.834
$ws.Range("J68").Value = 29929.834
$ws.Range("L68").Value = 29929.834
$ws.Range("N68").Value = -31427.834
$ws.Range("H70").Value = 37750
$ws.Range("J70").Value = 38000
$ws.Range("L70").Value = 38000
$ws.Range("N70").Value = -38630
$ws.Range("H71").Value = 29929.834
$ws.Range("J71").Value = 29929.834
$ws.Range("L71").Value = 89789.50199999999
$ws.Range("N71").Value = -97277.50199999999
$ws.Range("H73").Value = 37750
$ws.Range("J73").Value = 38000
$ws.Range("L73").Value = 38000
$ws.Range("N73").Value = -40184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 760.8837
$ws.Range("I5").Value = 459.43588
$ws.Range("K5").Value = 1378.30764
$ws.Range("M5").Value = -1266.30764
$ws.Range("H94").Value = 3739.3103
$ws.Range("J94").Value = 3940
$ws.Range("L94").Value = 11820
$ws.Range("N94").Value = -13172
$ws.Range("H121").Value = 9074.138999999999
$ws.Range("I121").Value = 212.85
$ws.Range("J121").Value = 20150.75
$ws.Range("K121").Value = 638.55
$ws.Range("L121").Value = 60452.25
$ws.Range("M121").Value = 671.45
$ws.Range("N121").Value = -63072.25
$ws.Range("H132").Value = 1764
$ws.Range("I132").Value = 1371.3636
$ws.Range("J132").Value = 1969.6666
$ws.Range("K132").Value = 12342.2724
$ws.Range("L132").Value = 17726.9994
$ws.Range("M132").Value = -9812.2724
$ws.Range("N132").Value = -22786.9994
$ws.Range("H135").Value = 760.8837
$ws.Range("I135").Value = 459.43588
$ws.Range("K135").Value = 4134.92292
$ws.Range("M135").Value = -1599.92292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 69183.164
$ws.Range("J52").Value = 80019.8
$ws.Range("L52").Value = 80019.8
$ws.Range("N52").Value = -80537.8
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 4430.7
$ws.Range("I122").Value = 3167.7878
$ws.Range("K122").Value = 9503.3634
$ws.Range("M122").Value = -7053.3634
$ws.Range("H125").Value = 22163
$ws.Range("J125").Value = 22163
$ws.Range("L125").Value = 22163
$ws.Range("N125").Value = -27083
$ws.Range("H132").Value = 3858.8684
$ws.Range("I132").Value = 3969.0476
$ws.Range("J132").Value = 3722.7646
$ws.Range("K132").Value = 11907.1428
$ws.Range("L132").Value = 11168.2938
$ws.Range("M132").Value = -9377.1428
$ws.Range("N132").Value = -16228.2938
$ws.Range("H140").Value = 26666.666
$ws.Range("J140").Value = 26666.666
$ws.Range("L140").Value = 26666.666
$ws.Range("N140").Value = -37026.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 167962.36
$ws.Range("I14").Value = 2500000
$ws.Range("J14").Value = 22210
$ws.Range("K14").Value = 2500000
$ws.Range("L14").Value = 22210
$ws.Range("M14").Value = -2499828
$ws.Range("N14").Value = -22554
$ws.Range("H17").Value = 23008.1
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H24").Value = 36000.332
$ws.Range("J24").Value = 36000.332
$ws.Range("L24").Value = 36000.332
$ws.Range("N24").Value = -36686.332
$ws.Range("H31").Value = 1838.9
$ws.Range("I31").Value = 842.9
$ws.Range("J31").Value = 2834.9
$ws.Range("K31").Value = 842.9
$ws.Range("L31").Value = 2834.9
$ws.Range("M31").Value = -594.9
$ws.Range("N31").Value = -3330.9
$ws.Range("H40").Value = 3014
$ws.Range("I40").Value = 2992.6155
$ws.Range("J40").Value = 3048.75
$ws.Range("K40").Value = 2992.6155
$ws.Range("L40").Value = 3048.75
$ws.Range("M40").Value = -2856.6155
$ws.Range("N40").Value = -3320.75
$ws.Range("H69").Value = 29975
$ws.Range("J69").Value = 29975
$ws.Range("L69").Value = 29975
$ws.Range("N69").Value = -31597
$ws.Range("H72").Value = 29975
$ws.Range("J72").Value = 29975
$ws.Range("L72").Value = 89925
$ws.Range("N72").Value = -98037
$ws.Range("H100").Value = 2116.0588
$ws.Range("I100").Value = 1429.125
$ws.Range("J100").Value = 2726.6667
$ws.Range("K100").Value = 1429.125
$ws.Range("L100").Value = 2726.6667
$ws.Range("M100").Value = -888.125
$ws.Range("N100").Value = -3808.6667
$ws.Range("H106").Value = 29000
$ws.Range("J106").Value = 29000
$ws.Range("L106").Value = 29000
$ws.Range("N106").Value = -31524
$ws.Range("H132").Value = 2500.0227
$ws.Range("I132").Value = 1762.04
$ws.Range("J132").Value = 3471.0527
$ws.Range("K132").Value = 5286.12
$ws.Range("L132").Value = 10413.1581
$ws.Range("M132").Value = -2756.12
$ws.Range("N132").Value = -15473.1581
$ws.Range("H136").Value = 4734.222
$ws.Range("I136").Value = 4019.9375
$ws.Range("J136").Value = 5773.1816
$ws.Range("K136").Value = 12059.8125
$ws.Range("L136").Value = 17319.5448
$ws.Range("M136").Value = -9509.8125
$ws.Range("N136").Value = -22419.5448
